$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 65
$ws.Range("F6").Value = 700
$ws.Range("F9").Value = 839
$ws.Range("F10").Value = 1555
$ws.Range("F11").Value = 1244
$ws.Range("E12").Value = "2024.06.22 10:00-06.22 16:20"
$ws.Range("F12").Value = 1489
$ws.Range("F13").Value = 47
$ws.Range("F14").Value = 1404
$ws.Range("F15").Value = 328
$ws.Range("F19").Value = 1070
$ws.Range("F20").Value = 354
$ws.Range("F23").Value = 1599
$ws.Range("F28").Value = 1165
$ws.Range("F29").Value = 288042
$ws.Range("F30").Value = 1019
$ws.Range("F34").Value = 1110
$ws.Range("F36").Value = 1106
$ws.Range("F37").Value = 58
$ws.Range("F38").Value = 232
$ws.Range("F41").Value = 1648
$ws.Range("F42").Value = 107
$ws.Range("F43").Value = 74
$ws.Range("F46").Value = 791

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 2562
$ws.Range("F10").Value = 1205
$ws.Range("F12").Value = 722
$ws.Range("F13").Value = 241
$ws.Range("F29").Value = 199
$ws.Range("F31").Value = 48
$ws.Range("F33").Value = 6
$ws.Range("F42").Value = 61

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 255
$ws.Range("F5").Value = 2861
$ws.Range("F6").Value = 4601
$ws.Range("F9").Value = 565
$ws.Range("F10").Value = 706
$ws.Range("F11").Value = 452
$ws.Range("F12").Value = 310
$ws.Range("F13").Value = 977
$ws.Range("F14").Value = 267
$ws.Range("F15").Value = 606

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 255
$ws.Range("F4").Value = 4601
$ws.Range("F5").Value = 706
$ws.Range("F6").Value = 65
$ws.Range("F7").Value = 310
$ws.Range("F8").Value = 310
$ws.Range("F9").Value = 977
$ws.Range("F10").Value = 977
$ws.Range("F13").Value = 839
$ws.Range("F14").Value = 2562
$ws.Range("F15").Value = 1205
$ws.Range("F16").Value = 1555
$ws.Range("F17").Value = 1244
$ws.Range("E18").Value = "2024.06.22 10:00-06.22 16:20"
$ws.Range("F18").Value = 1489
$ws.Range("F19").Value = 47
$ws.Range("F20").Value = 1404
$ws.Range("F21").Value = 241
$ws.Range("G21").Value = 180
$ws.Range("F22").Value = 328
$ws.Range("F24").Value = 1650
$ws.Range("F26").Value = 1070
$ws.Range("F27").Value = 354
$ws.Range("F28").Value = 606
$ws.Range("F29").Value = 606
$ws.Range("F31").Value = 1599
$ws.Range("F36").Value = 1165
$ws.Range("F38").Value = 1019
$ws.Range("F40").Value = 1110
$ws.Range("F42").Value = 1106
$ws.Range("F44").Value = 232
$ws.Range("F47").Value = 1648
$ws.Range("F48").Value = 107
$ws.Range("F49").Value = 74
$ws.Range("F52").Value = 791
